$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows above row 599 (everything from old row 599 onward shifts down by 5)
$ws.Rows("599:603").Insert()

# Shared/common column values for this vendor's Tomate rows
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$catId     = 100112020
$categoria = "Tomate"
$unidad    = "`$/bandeja 18 kilos"
$origen    = "Provincia de Limarí"
$kgUnid    = 18
$clasif    = "Hortaliza"

$fecha = 44714

$newRows = @(
    @{ Row=599; Variedad="Larga vida"; Calidad="Primera"; Volumen=900;  PrecioMin=11000; PrecioMax=12000; PrecioProm=11500; PrecioKg=639 },
    @{ Row=600; Variedad="Larga vida"; Calidad="Segunda"; Volumen=800;  PrecioMin=9000;  PrecioMax=10000; PrecioProm=9500;  PrecioKg=528 },
    @{ Row=601; Variedad="Larga vida"; Calidad="Tercera"; Volumen=600;  PrecioMin=7000;  PrecioMax=8000;  PrecioProm=7500;  PrecioKg=417 },
    @{ Row=602; Variedad="Semiduro";   Calidad="Primera"; Volumen=400;  PrecioMin=5000;  PrecioMax=6000;  PrecioProm=5500;  PrecioKg=306 },
    @{ Row=603; Variedad="Semiduro";   Calidad="Segunda"; Volumen=400;  PrecioMin=3000;  PrecioMax=4000;  PrecioProm=3500;  PrecioKg=194 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PrecioMin
    $ws.Cells.Item($row, 12).Value = $r.PrecioMax
    $ws.Cells.Item($row, 13).Value = $r.PrecioProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PrecioKg
    $ws.Cells.Item($row, 17).Value = $kgUnid
    $ws.Cells.Item($row, 18).Value = $clasif
}
